# This script reproduces the automated "Updated cryptos list ... with GitHub
# Actions" commit: it refreshes the Price (D) and Volume(1h) (E) columns for
# the existing 50 ranked coins, and also fixes rows 9-10 (USDC and Cardano had
# been swapped) and row 51 (now tracks Aave instead of Cosmos).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D ("Price") stores plain numeric-looking text (e.g. "1.00", "0.0000284").
# Plain Range.Value assignment would make Excel "helpfully" reinterpret those as
# real numbers (dropping trailing zeros, collapsing "0.0000200" to scientific, etc),
# so pre-format the affected cells as Text first to preserve the exact source text.
$ws.Range('D4:D10').NumberFormat = "@"
$ws.Range('D12:D15').NumberFormat = "@"
$ws.Range('D18:D19').NumberFormat = "@"
$ws.Range('D21:D28').NumberFormat = "@"
$ws.Range('D30:D32').NumberFormat = "@"
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D36:D38').NumberFormat = "@"
$ws.Range('D40:D43').NumberFormat = "@"
$ws.Range('D45:D51').NumberFormat = "@"

$ws.Range('D2').Value = '95.000.10'
$ws.Range('E2').Value = '  -1.91%  '

$ws.Range('D3').Value = '3.578.50'
$ws.Range('E3').Value = '  -2.84%  '

$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  +0.16%  '

$ws.Range('D5').Value = '2.26'
$ws.Range('E5').Value = '  +17.18%  '

$ws.Range('D6').Value = '223.64'
$ws.Range('E6').Value = '  -5.34%  '

$ws.Range('D7').Value = '632.00'
$ws.Range('E7').Value = '  -3.87%  '

$ws.Range('D8').Value = '0.406'
$ws.Range('E8').Value = '  -4.05%  '

$ws.Range('B9').Value = 'USDC'
$ws.Range('C9').Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Range('D9').Value = '1.00'
$ws.Range('E9').Value = '  +0.14%  '

$ws.Range('B10').Value = 'Cardano'
$ws.Range('C10').Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range('D10').Value = '1.06'
$ws.Range('E10').Value = '  -0.47%  '

$ws.Range('D11').Value = '3.577.24'
$ws.Range('E11').Value = '  -2.80%  '

$ws.Range('D12').Value = '45.63'
$ws.Range('E12').Value = '  +3.55%  '

$ws.Range('D13').Value = '0.204'
$ws.Range('E13').Value = '  -2.00%  '

$ws.Range('D14').Value = '0.0000284'
$ws.Range('E14').Value = '  -5.23%  '

$ws.Range('D15').Value = '6.39'
$ws.Range('E15').Value = '  -5.16%  '

$ws.Range('D16').Value = '4.249.87'
$ws.Range('E16').Value = '  -2.78%  '

$ws.Range('D17').Value = '94.729.87'
$ws.Range('E17').Value = '  -2.05%  '

$ws.Range('D18').Value = '8.71'
$ws.Range('E18').Value = '  -3.71%  '

$ws.Range('D19').Value = '19.77'
$ws.Range('E19').Value = '  +5.83%  '

$ws.Range('D20').Value = '3.583.06'
$ws.Range('E20').Value = '  -2.92%  '

$ws.Range('D21').Value = '12.74'
$ws.Range('E21').Value = '  -1.67%  '

$ws.Range('D22').Value = '0.504'
$ws.Range('E22').Value = '  -0.91%  '

$ws.Range('D23').Value = '496.94'
$ws.Range('E23').Value = '  -4.71%  '

$ws.Range('D24').Value = '3.19'
$ws.Range('E24').Value = '  -5.85%  '

$ws.Range('D25').Value = '0.236'
$ws.Range('E25').Value = '  +14.68%  '

$ws.Range('D26').Value = '116.02'
$ws.Range('E26').Value = '  +14.84%  '

$ws.Range('D27').Value = '0.0000200'
$ws.Range('E27').Value = '  -4.97%  '

$ws.Range('D28').Value = '6.68'
$ws.Range('E28').Value = '  -3.12%  '

$ws.Range('D29').Value = '3.781.24'
$ws.Range('E29').Value = '  -2.64%  '

$ws.Range('D30').Value = '12.40'
$ws.Range('E30').Value = '  -6.79%  '

$ws.Range('D31').Value = '12.59'
$ws.Range('E31').Value = '  +1.43%  '

$ws.Range('D32').Value = '2.86'
$ws.Range('E32').Value = '  -4.68%  '

$ws.Range('E33').Value = '  -0.07%  '

$ws.Range('D34').Value = '0.999'
$ws.Range('E34').Value = '  +0.02%  '

$ws.Range('E35').Value = '  -5.72%  '

$ws.Range('D36').Value = '1.74'
$ws.Range('E36').Value = '  -6.15%  '

$ws.Range('D37').Value = '31.27'
$ws.Range('E37').Value = '  -2.66%  '

$ws.Range('D38').Value = '0.575'
$ws.Range('E38').Value = '  -2.09%  '

$ws.Range('E39').Value = '  -0.02%  '

$ws.Range('D40').Value = '583.86'
$ws.Range('E40').Value = '  -9.29%  '

$ws.Range('D41').Value = '8.19'
$ws.Range('E41').Value = '  -6.77%  '

$ws.Range('D42').Value = '6.69'
$ws.Range('E42').Value = '  -2.02%  '

$ws.Range('D43').Value = '40.33'
$ws.Range('E43').Value = '  -0.48%  '

$ws.Range('E44').Value = '  -1.82%  '

$ws.Range('D45').Value = '0.462'
$ws.Range('E45').Value = '  -0.85%  '

$ws.Range('D46').Value = '0.0464'
$ws.Range('E46').Value = '  +1.53%  '

$ws.Range('D47').Value = '1.88'
$ws.Range('E47').Value = '  -7.70%  '

$ws.Range('D48').Value = '0.906'
$ws.Range('E48').Value = '  -4.71%  '

$ws.Range('D49').Value = '23.42'
$ws.Range('E49').Value = '  -0.82%  '

$ws.Range('D50').Value = '3.60'
$ws.Range('E50').Value = '  +2.48%  '

$ws.Range('B51').Value = 'Aave'
$ws.Range('C51').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D51').Value = '216.35'
$ws.Range('E51').Value = '  +6.03%  '
